$wb = $excel.ActiveWorkbook

# Rename the "language_English" sheet to "language_ENG"
$langSheet = $wb.Worksheets.Item("language_English")
$langSheet.Name = "language_ENG"

# Update the "settings" sheet: the "default language name" row (A2) changes
# from "English" to "ENG"
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "ENG"
[void]$settings.Range("A2").Select()

# Restore the "tags" sheet as the active/selected tab
$tags = $wb.Worksheets.Item("tags")
[void]$tags.Activate()
